# [RCE] Update data for 585000 rows of data
#
# The source "ExperimentRecord" sheet tracks one row of metrics per
# training-set size. A new size (585000 = 520000 + 65000, i.e. the next
# step in the existing 65000-row stride) was appended as a new row
# directly under the previous last data row (row 17), pushing the blank
# spacer row and the trailing "Notes" legend rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 18 (shifts old row 18 -> 19, ..., old row 25 -> 26)
$ws.Rows.Item(18).Insert()

# Bring over row 17's formatting (borders/number-formats) so the new
# row looks like the rest of the data block instead of like the
# now-shifted blank spacer row.
$ws.Range("A17:N17").Copy()
$ws.Range("A18:N18").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row 18 with the next experiment's data.
$ws.Range("A18").Formula = "=A17+65000"
$ws.Range("B18").Formula = "=A18-D18"
$ws.Range("C18").Formula = "=B18/A18"
$ws.Range("D18").Value = 416849
$ws.Range("E18").Value = 0.903
$ws.Range("F18").Value = 0.922
$ws.Range("G18").Value = 0.923
$ws.Range("H18").Value = 0.525
$ws.Range("I18").Value = 0.908
$ws.Range("J18").Value = 0.869
$ws.Range("K18").Formula = "=2024.97/60"
$ws.Range("L18").Formula = "=61748.06/60"
$ws.Range("M18").Formula = "=60.24/60"
$ws.Range("N18").Value = 131

# Match the author's final cursor position/selection.
$ws.Range("A19").Select() | Out-Null
